$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the refreshed coin data (price + 1h volume change columns).
# Values that would otherwise be auto-detected by Excel as numbers get a
# leading apostrophe so they are stored as literal text, matching the
# original inline-string cell type used throughout this sheet.
$ws.Range("D2").Value = '26.138.49'
$ws.Range("E2").Value = '  -0.27%  '
$ws.Range("D3").Value = '1.654.82'
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").Value = '''218.56'
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("D6").Value = '''0.5313'
$ws.Range("E6").Value = '  +1.71%  '
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("D8").Value = '''0.2613'
$ws.Range("E8").Value = '  -2.07%  '
$ws.Range("D9").Value = '''0.06339'
$ws.Range("E9").Value = '  +0.27%  '
$ws.Range("D10").Value = '''20.45'
$ws.Range("E10").Value = '  -2.95%  '
$ws.Range("D11").Value = '''0.07749'
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '''4.494'
$ws.Range("E12").Value = '  +1.51%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.639.67'
$ws.Range("E13").Value = '  -1.18%  '
$ws.Range("D14").Value = '''0.5472'
$ws.Range("E14").Value = '  +0.22%  '
$ws.Range("D15").Value = '0.0₅8139'
$ws.Range("E15").Value = '  -1.07%  '
$ws.Range("D16").Value = '''65.35'
$ws.Range("E16").Value = '  +0.80%  '
$ws.Range("D17").Value = '26.145.28'
$ws.Range("E17").Value = '  -0.40%  '
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").Value = '''4.553'
$ws.Range("E19").Value = '  -2.29%  '
$ws.Range("D20").Value = '''193.87'
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("D21").Value = '''10.06'
$ws.Range("E21").Value = '  -0.71%  '
$ws.Range("D22").Value = '''6.005'
$ws.Range("E22").Value = '  -1.08%  '
$ws.Range("D23").Value = '''1.003'
$ws.Range("E23").Value = '  -0.36%  '
$ws.Range("D24").Value = '''140.38'
$ws.Range("E24").Value = '  +1.24%  '
$ws.Range("D25").Value = '''0.1243'
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("D26").Value = '''7.280'
$ws.Range("E26").Value = '  +0.81%  '
$ws.Range("D27").Value = '''16.21'
$ws.Range("E27").Value = '  +0.54%  '
$ws.Range("E28").Value = '  +1.77%  '
$ws.Range("D29").Value = '''0.05948'
$ws.Range("E29").Value = '  -0.74%  '
$ws.Range("D30").Value = '''1.280'
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("D31").Value = '''3.512'
$ws.Range("E31").Value = '  -3.96%  '
$ws.Range("E32").Value = '  -2.15%  '
$ws.Range("D33").Value = '''1.547'
$ws.Range("E33").Value = '  -5.19%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '''0.9485'
$ws.Range("E34").Value = '  -3.10%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''2.411'
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("D36").Value = '''2.760'
$ws.Range("E36").Value = '  -0.82%  '
$ws.Range("D37").Value = '''0.5650'
$ws.Range("E37").Value = '  -4.02%  '
$ws.Range("D38").Value = '''0.01613'
$ws.Range("E38").Value = '  +1.24%  '
$ws.Range("D39").Value = '''5.861'
$ws.Range("E39").Value = '  -1.41%  '
$ws.Range("D40").Value = '''0.8478'
$ws.Range("E40").Value = '  -1.37%  '
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '''101.08'
$ws.Range("E42").Value = '  +1.55%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '1.010.70'
$ws.Range("E43").Value = '  -2.01%  '
$ws.Range("D44").Value = '1.800.26'
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").Value = '''56.88'
$ws.Range("E45").Value = '  -0.48%  '
$ws.Range("D46").Value = '0.0₈106'
$ws.Range("E46").Value = '  -6.37%  '
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("D48").Value = '''0.4287'
$ws.Range("E48").Value = '  +1.38%  '
$ws.Range("D49").Value = '''1.478'
$ws.Range("E49").Value = '  +0.59%  '
$ws.Range("E50").Value = '  -0.67%  '
$ws.Range("D51").Value = '''7.734'
$ws.Range("E51").Value = '  -4.50%  '
